# Insert two new weekly price rows for "Poroto verde" (Comercializadora del
# Agro de Limarí) right after the existing row 144, pushing the previous
# rows 145:174 down to 147:176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 145:174 down by two to make room for the new records.
$ws.Rows("145:146").Insert()

# Common (unchanged-pattern) column values shared by every data row in this
# block.
$mercado = 2
$comercializadora = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112031
$categoria = "Poroto verde"
$calidad = "Primera"
$origen = "Provincia de Limarí"
$kgUnidades = 25
$clasificacion = "Hortaliza"

function Set-PrecioRow($Row, $Fecha, $Variedad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Unidad, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $mercado
    $ws.Cells.Item($Row, 2).Value = $comercializadora
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $categoriaId
    $ws.Cells.Item($Row, 7).Value = $categoria
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $kgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasificacion
}

Set-PrecioRow 145 44637 "Magnum" 700 19000 21000 20000 "$/malla 25 kilos" 800
Set-PrecioRow 146 44637 "Sin especificar" 440 23000 25000 24000 "$/malla 25 kilos" 960
